$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.093599557876587
$ws.Range("B1").Value = 1.105777025222778
$ws.Range("C1").Value = 1.922296404838562
$ws.Range("D1").Value = 3.694671630859375
$ws.Range("E1").Value = 5.68695068359375
